$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the "day of the week"
#    paragraph to the end of the "Determinar cuanto le cuesta..." paragraph.
# ---------------------------------------------------------------------------

# Remove the existing (hidden) _GoBack bookmark.
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

# Locate the target paragraph.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    $txt = $para.Range.Text
    if ($txt -like "*Determinar cu*nto le cuesta a una persona contratar una p*liza*") {
        $targetPara = $para
        break
    }
}

# Insert a temporary one-character marker right after the paragraph's text
# (but before its paragraph mark), wrap a fresh "_GoBack" bookmark tightly
# around that single character, then delete the marker again. This leaves
# the bookmark collapsed exactly between the run and the paragraph mark --
# i.e. right where the original bookmark used to sit in the other paragraph.
$pr = $targetPara.Range
$pr.InsertAfter("Z")

$newEnd = $targetPara.Range.End
$markerStart = $newEnd - 2
$markerEnd = $newEnd - 1
$markerRange = $d.Range($markerStart, $markerEnd)
$d.Bookmarks.Add("_GoBack", $markerRange)

$delRange = $d.Range($markerStart, $markerEnd)
$delRange.Text = ""

# ---------------------------------------------------------------------------
# 2) Color the "calificacion" paragraph's text (and paragraph mark) red.
# ---------------------------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    $txt = $para.Range.Text
    if ($txt -like "*Realizar algoritmo que, con base en una calificaci*n proporcionada*") {
        $para.Range.Font.Color = 255
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Merge the three runs (split by a spell-check proofErr wrapper around
#    "mas") describing the Feb-14th gift scenario back into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "El 14 de febrero una persona desea comprarle al ser querido que mas aprecia en ese momento, su dilema radica en regalo puede hacerle, las alternativas que tienen son los siguientes:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El 14 de febrero una persona desea comprarle al ser querido que mas aprecia en ese momento, su dilema radica en regalo puede hacerle, las alternativas que tienen son los siguientes:",
    2) | Out-Null
